$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Detect_Actual_Distance (column R) values
$ws.Range("R2").Value = 21.62267899129298
$ws.Range("R3").Value = 19.69708450735381
$ws.Range("R5").Value = 22.94250238250783
$ws.Range("R7").Value = 56.05221049722274
$ws.Range("R8").Value = 54.54562755004776

# Row 11: update actual lane/edge identifiers and distance to reflect
# correctly resolved lane geometry (Campus_WB instead of gneE1.93)
$ws.Range("N11").Value = "Campus_WB_1"
$ws.Range("P11").Value = "Campus_WB"
$ws.Range("R11").Value = 13.28494941060483
